$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: rename monthly_payment -> weekly_payment (column M) ---
$ws.Range("M1").Value = "weekly_payment"

# --- Clear the old single-customer row (row 2) ---
# The row is kept in place (now blank) while the real data moves into the
# new multi-loan block below it.
for ($col = 1; $col -le 20; $col++) {
    $cell = $ws.Cells.Item(2, $col)
    $cell.Value = ""
    # touch a formatting property (no-op) so the now-empty cell still gets
    # written out/kept instead of being dropped entirely
    $bold = $cell.Font.Bold
    $cell.Font.Bold = $bold
}

# --- New customer / loan rows starting at row 3 ---
$rows = @(
    @("8325","CHARNEQUA BENNET",16990,1440.9,820,0,799,1185.21,3000,0,238,28.054,136.23,7,"`$18,235.11","02/24/1983","595-20-8868","charnequabennet02@gmail.com","(772) 301-2514","745 North Powerline Road, Pompano Beach FL 33069"),
    @("8326","CHAD WHYMS",11500,1190.7,820,0,799,855.8100000000001,2500,0,208,28.554,102.3,7,"`$12,665.51","03/17/1984","593-64-4629","chadwyms@yahoo.com","(786) 617-2472","3595 NW 187TH ST, MIAMI GARDENS FL 33056"),
    @("8322","JEANNE LUBIN",11300,1128.9,820,0,799,843.8100000000001,2000,0,208,28.556,104.13,7,"`$12,891.71","01/20/1991","766-84-9187","jeannelubin1991@gmail.com","(786) 805-8896","820 Northwest 81st Street, Miami FL 33150"),
    @("8328","DIANNY MEDINA",11000,1438.6,820,0,799,825.8100000000001,2200,0,238,28.05,94.76000000000001,7,"`$12,683.41","07/30/1970","879-30-9526","DIANNYMEDINA500@GMAIL.COM","(786) 231-7929","16720 Northwest 48th Court, Miami Gardens FL 33055"),
    @("8305","SHAKIERA THOMPSON",21990,1235.7,820,0,799,1485.21,4500,0,238,28.05,163.09,7,"`$21,829.91","03/21/1989","591-80-0690","SHAKIERA21@ICLOUD.COM","(229) 343-9706","3960 55TH ST NORTH  APT 101, ST PETERSBURG FL 33709")
)

# columns that must stay as literal text even though they look numeric /
# date-like (loan number, computed "amount financed" string, DOB string)
$textCols = @(1, 15, 16)

$r = 3
foreach ($row in $rows) {
    for ($col = 1; $col -le 20; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $value = $row[$col - 1]
        if ($textCols -contains $col) {
            $cell.Value = "'" + $value
        } else {
            $cell.Value = $value
        }
    }
    $r++
}

Write-Output "done"
